$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update column B (prediction distance/score) for rows 2..108 ---
$bValues = @(2450.3883699564058,2696.8669210263397,2359.8443429788226,2609.3986497324172,1859.6607549073278,3638.5791972500438,1862.3095663409281,2886.0783273804786,2767.8482318655815,3069.7745961717678,3625.5374398950953,2337.4229952855994,2333.2105051907738,2205.1740855836206,2013.075379025393,1803.0745621228889,2494.3826835488394,1617.0076792475311,1993.3920930649942,2373.4904017205849,2482.9601244215637,2151.3147821457032,2526.4258398085799,3085.3121084004847,2700.9580480812579,2231.764019114134,2200.5349373358149,1841.4077973896397,2092.9721007430962,2116.961804656934,2567.5304221340675,2543.0436961362921,2720.3230766338697,2326.1549804271644,2343.033841259999,2458.9154794080746,2497.2336538096606,1769.1186651795315,2268.1103493241849,2066.9042173109142,2114.2210655248082,3095.6309603517425,2048.6468072794314,2283.8572445444879,3842.0342369756318,2426.5329916867477,2091.5985352889456,2219.3276546927314,2538.9585492323513,2436.076886257461,1872.3640243249411,2015.3443634012765,2617.7695907237107,2422.7299162800023,2362.2620082649823,2442.4246025622815,2773.193728069602,1355.2251873695523,2835.0683694009836,2847.9942108634582,1655.9057895821934,3075.9989372366581,3361.6054234392568,2472.2045457962322,3253.9803857732541,3229.7036797482442,1902.8451235938942,3118.9321109000789,962.32307369096281,1517.0226512506019,3551.2287674687605,3481.7158850944397,3452.2834039909631,3691.4748574407613,2127.3889853216324,3203.1147411535972,2016.6463498447033,2346.934719291185,3385.954069690295,3343.2199636648584,2232.6506532569711,2799.4406481073224,1772.3276428509344,3864.5907157648126,2403.7641461126632,2958.1241943757573,2202.5239542653867,2065.5503605847402,1592.4097376315494,3177.1809811016337,2894.5625526247973,2678.2892768783959,2753.5290970539772,3534.2147963383368,3364.4942886610952,3380.5453985933118,2989.7437351939193,3427.8575475769048,2750.1395262257993,3046.6927302833351,2401.57166098677,2401.57166098677,1278.6723586869375,2186.4123533004408,3312.198688585624,2545.221903134403,2545.221903134403)
for ($i = 0; $i -lt $bValues.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value2 = $bValues[$i]
}

# --- Update column A (genome id) for rows whose referenced shared-string shifted ---
# (genome "even_MAG-GUT81204.fa" and "even_MAG-GUT7042.fa" were moved to the end of the list)
$aRows  = @(94,95,96,97,98,99,100,101,102,103,104,105,106,107,108)
$aTexts = @("even_MAG-GUT74347.fa","even_MAG-GUT7772.fa","even_MAG-GUT80568.fa","even_MAG-GUT85070.fa","even_MAG-GUT86606.fa","even_MAG-GUT88052.fa","even_MAG-GUT88444.fa","even_MAG-GUT9016.fa","even_MAG-GUT91733.fa","even_MAG-GUT91735.fa","even_MAG-GUT91947.fa","even_MAG-GUT92065.fa","even_MAG-GUT9523.fa","even_MAG-GUT7042.fa","even_MAG-GUT81204.fa")
for ($i = 0; $i -lt $aRows.Count; $i++) {
    $ws.Cells.Item($aRows[$i], 1).Value2 = $aTexts[$i]
}
